$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.054.95"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "2.301.38"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("Z1").Formula = "=""299.94"""
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("Z1").Formula = "=""97.67"""
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("Z1").Formula = "=""33.83"""
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("Z1").Formula = "=""49.10"""
$ws.Range("Z1").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("Z1").Formula = "=""17.15"""
$ws.Range("Z1").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E14").Value = "  +11.25%  "
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "2.659.59"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "2.345.26"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").Value = "42.983.53"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("Z1").Formula = "=""11.64"""
$ws.Range("Z1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("Z1").Formula = "=""236.67"""
$ws.Range("Z1").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  +5.42%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("Z1").Formula = "=""2.45"""
$ws.Range("Z1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("Z1").Formula = "=""24.38"""
$ws.Range("Z1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("Z1").Formula = "=""166.28"""
$ws.Range("Z1").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("Z1").Formula = "=""33.84"""
$ws.Range("Z1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("Z1").Formula = "=""4.60"""
$ws.Range("Z1").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E35").Value = "  +6.37%  "
$ws.Range("Z1").Formula = "=""2.43"""
$ws.Range("Z1").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E36").Value = "  +1.79%  "
$ws.Range("Z1").Formula = "=""16.89"""
$ws.Range("Z1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E37").Value = "  +4.33%  "
$ws.Range("Z1").Formula = "=""0.0702"""
$ws.Range("Z1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("D44").Value = "1.995.44"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("Z1").Formula = "=""0.0285"""
$ws.Range("Z1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("Z1").Formula = "=""9.82"""
$ws.Range("Z1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("Z1").Formula = "=""17.53"""
$ws.Range("Z1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "2.525.21"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("Z1").Formula = "=""53.40"""
$ws.Range("Z1").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("E51").Value = "  -1.77%  "
